$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '81.553.97'
$ws.Range("E2").Value = '  +2.72%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.166.97'
$ws.Range("E3").Value = '  -0.92%  '

$ws.Range("E4").Value = '  +0.01%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '210.46'
$ws.Range("E5").Value = '  +2.24%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '621.06'
$ws.Range("E6").Value = '  -2.17%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.278'
$ws.Range("E7").Value = '  +18.37%  '

$ws.Range("E8").Value = '  -0.06%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.581'
$ws.Range("E9").Value = '  -0.76%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '3.163.39'
$ws.Range("E10").Value = '  -0.92%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.580'
$ws.Range("E11").Value = '  -0.64%  '

$ws.Range("E13").Value = '  -0.36%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.29'
$ws.Range("E14").Value = '  -4.82%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.742.54'
$ws.Range("E15").Value = '  -0.98%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '31.33'
$ws.Range("E16").Value = '  -1.22%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '81.284.29'
$ws.Range("E17").Value = '  +2.56%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.154.60'
$ws.Range("E18").Value = '  -0.91%  '

$ws.Range("E19").Value = '  +0.45%  '

$ws.Range("E20").Value = '  -4.60%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '430.25'
$ws.Range("E21").Value = '  +0.27%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '8.94'
$ws.Range("E22").Value = '  -2.23%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.06'
$ws.Range("E23").Value = '  +1.27%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '7.23'
$ws.Range("E24").Value = '  +5.25%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '5.19'
$ws.Range("E25").Value = '  +8.36%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.319.81'
$ws.Range("E26").Value = '  -1.16%  '

$ws.Range("E27").Value = '  -0.68%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.73'
$ws.Range("E28").Value = '  -4.52%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.998'
$ws.Range("E29").Value = '  -0.23%  '

$ws.Range("E30").Value = '  +1.76%  '

$ws.Range("B31").Value = 'Bittensor'
$ws.Range("C31").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '579.60'
$ws.Range("E31").Value = '  +9.86%  '

$ws.Range("B32").Value = 'Binance-PegBSC-USD'
$ws.Range("C32").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.995'
$ws.Range("E32").Value = '  -0.53%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '8.89'
$ws.Range("E33").Value = '  -1.07%  '

$ws.Range("E34").Value = '  +0.02%  '

$ws.Range("E35").Value = '  +7.34%  '

$ws.Range("B36").Value = 'Cronos'
$ws.Range("C36").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.139'
$ws.Range("E36").Value = '  +15.93%  '

$ws.Range("B37").Value = 'PancakeSwap'
$ws.Range("C37").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.98'
$ws.Range("E37").Value = '  -1.31%  '

$ws.Range("E38").Value = '  -1.15%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.999'
$ws.Range("E39").Value = '  -0.02%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.06'
$ws.Range("E40").Value = '  +10.74%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.406'
$ws.Range("E41").Value = '  +0.21%  '

$ws.Range("E42").Value = '  +13.48%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '20.75'
$ws.Range("E43").Value = '  +3.67%  '

$ws.Range("E44").Value = '  +18.42%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '159.35'
$ws.Range("E45").Value = '  -3.19%  '

$ws.Range("E46").Value = '  +0.01%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '186.75'
$ws.Range("E47").Value = '  -3.35%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '45.06'
$ws.Range("E48").Value = '  +4.47%  '

$ws.Range("E49").Value = '  +0.45%  '

$ws.Range("E50").Value = '  -4.74%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '25.78'
$ws.Range("E51").Value = '  -0.48%  '

